$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price reading was added for "Feria Lagunitas de Puerto Montt -
# Espinaca". It belongs chronologically as the most recent entry, so it is
# inserted as the new row 18 and every existing record from the old row 18
# down to the old row 49 shifts down by one row (old 49 becomes new 50).
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 4
$ws.Range("B18").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C18").Value = 'Los Lagos'
$ws.Range("D18").Value = '2022-10-04'
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 100112012
$ws.Range("G18").Value = 'Espinaca'
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 12000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 12000
$ws.Range("N18").Value = '$/cuna 10 kilos'
$ws.Range("O18").Value = 'Región Metropolitana'
$ws.Range("P18").Value = 1200
$ws.Range("Q18").Value = 10
$ws.Range("R18").Value = 'Hortaliza'
